$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (columns A-F): Giorno, TISG_PDR_G, fcs, buy_BEE_MWH, sell_lago_MWH, need_to_buy_MW
$data = @(
    @(2,  45870, 5173.64085737321,  4183.23451144358,  1944, 5457.066558, 105.110842169599),
    @(3,  45871, 929.483082685461,  1590.01054658509,  1944, 1755.527894, 19.6689732458178),
    @(4,  45872, 821.86701844334,   1556.56142387528,  1944, 1638.189309, 17.8701547679975),
    @(5,  45873, 3852.00481440493, 3775.12601353845,  1944, 4263.232285, 93.43139517223),
    @(6,  45874, 3852.00481440493, 3739.48013184365,  1944, 4263.232285, 91.9461501016136),
    @(7,  45875, 3939.76634278604, 3723.23814077569,  1944, 4337.311238, 90.6992931662351),
    @(8,  45876, 3939.76634278604, 3692.32232372257,  1944, 4337.311238, 89.4111341223554),
    @(9,  45877, 3939.76634278604, 3634.22042357079,  1944, 4337.311238, 86.9902216160313),
    @(10, 45878, 715.68594436081,  1252.3999908627,   1944, 1483.420019, 3.17225272924541),
    @(11, 45879, 634.559148236726, 1248.97653357819,  1944, 1394.227706, 2.69354547256082),
    @(12, 45880, 3735.57396241807, 3503.37258922662,  1944, 4141.254797, 81.8772259920228),
    @(13, 45881, 3735.57396241807, 3493.93914833543,  1944, 4141.254797, 81.4841659548898),
    @(14, 45882, 3735.57396241807, 3469.96937540677,  1944, 4141.254797, 80.4854254161957),
    @(15, 45883, 3735.57396241807, 3461.18696271948,  1944, 4141.254797, 80.1194915542252)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
